# Apply the "model_performance" presentation-prep edits:
#  - fix "Random Forst" typo -> "Random Forest"
#  - shade the header row with a light gray fill, center the metric headers
#  - widen column A to fit the corrected text
#  - zoom the sheet in to 250% and move the active selection to D12

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the shared string "Random Forst" (row 5, column A).
$ws.Range("A5").Value = "Random Forest"

# Header row formatting: bold text on a light gray fill; center the metric
# columns (Accuracy/Recall/Precision/F1) while the "Model" header stays
# left-aligned.
$headerRange = $ws.Range("A1:E1")
$headerRange.Interior.ColorIndex = 15
$headerRange.Font.Bold = $true

$ws.Range("B1:E1").HorizontalAlignment = -4108

# Widen column A to fit "Random Forest" / the other labels.
$ws.Columns.Item(1).ColumnWidth = 27.1640625

# Zoom in and move the selection.
$ws.Application.ActiveWindow.Zoom = 250
$ws.Range("D12").Select()
